$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns B:E in the data range to Text format first so that
# numeric-looking strings (e.g. "1.005", "45.50") are preserved exactly
# as typed instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.199.04"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.792.20"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "337.67"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").Value = "0.4522"
$ws.Range("E7").Value = "  +20.03%  "

$ws.Range("D8").Value = "0.3576"
$ws.Range("E8").Value = "  +6.53%  "

$ws.Range("D9").Value = "45.50"
$ws.Range("E9").Value = "  -0.65%  "

$ws.Range("D10").Value = "1.139"
$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("D11").Value = "0.07463"
$ws.Range("E11").Value = "  +3.69%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").Value = "22.37"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "6.208"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").Value = "7.230"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").Value = "1.791.86"
$ws.Range("E16").Value = "  +1.71%  "

$ws.Range("D17").Value = "0.00001083"
$ws.Range("E17").Value = "  +2.82%  "

$ws.Range("E18").Value = "  +1.83%  "

$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "17.18"
$ws.Range("E21").Value = "  +1.12%  "

$ws.Range("D22").Value = "6.383"
$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").Value = "28.209.00"
$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").Value = "11.83"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").Value = "2.380"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "20.42"
$ws.Range("E26").Value = "  +3.32%  "

$ws.Range("D27").Value = "153.65"
$ws.Range("E27").Value = "  +0.38%  "

$ws.Range("D28").Value = "2.371"
$ws.Range("E28").Value = "  +0.53%  "

$ws.Range("D29").Value = "1.995.53"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "132.28"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.265"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "4.071"
$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("D33").Value = "5.862"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "0.09398"
$ws.Range("E34").Value = "  +7.57%  "

$ws.Range("D35").Value = "0.02364"
$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").Value = "12.07"
$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("D37").Value = "0.6648"
$ws.Range("E37").Value = "  +0.56%  "

$ws.Range("D38").Value = "0.06222"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "0.2155"
$ws.Range("E39").Value = "  +2.16%  "

$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("D41").Value = "1.482"
$ws.Range("E41").Value = "  +2.52%  "

$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").Value = "8.044"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "0.9992"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").Value = "13.94"
$ws.Range("E45").Value = "  +1.93%  "

$ws.Range("D46").Value = "3.857"
$ws.Range("E46").Value = "  +0.49%  "

$ws.Range("D47").Value = "0.6055"
$ws.Range("E47").Value = "  +0.18%  "

$ws.Range("D48").Value = "128.11"
$ws.Range("E48").Value = "  -1.33%  "

$ws.Range("D49").Value = "2.019"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("E50").Value = "  -2.23%  "

$ws.Range("D51").Value = "1.160"
$ws.Range("E51").Value = "  -1.34%  "

# Restore the original "General" number format so the cells' display
# format matches the source workbook (only their text content changed).
$dataRange.NumberFormat = "General"
